$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.772.60"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "'2.322.38"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'302.24"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "'94.08"
$ws.Range("E6").Value = "  -3.65%  "
$ws.Range("D7").Value = "'0.500"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("D9").Value = "'0.491"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("D10").Value = "'33.87"
$ws.Range("E10").Value = "  -4.66%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0779"
$ws.Range("E11").Value = "  -2.45%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "'18.64"
$ws.Range("E12").Value = "  -4.87%  "
$ws.Range("D13").Value = "'0.121"
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").Value = "'6.69"
$ws.Range("E14").Value = "  -3.70%  "
$ws.Range("D15").Value = "'2.687.24"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "'2.365.46"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").Value = "'0.786"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "'42.715.03"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "'11.96"
$ws.Range("E19").Value = "  -5.42%  "
$ws.Range("D20").Value = "'6.18"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").Value = "'0.0₃0884"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").Value = "'67.75"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "'234.77"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "'2.22"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("D27").Value = "'24.45"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("D28").Value = "'2.35"
$ws.Range("E28").Value = "  +14.06%  "
$ws.Range("D29").Value = "'9.09"
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("D30").Value = "'31.11"
$ws.Range("E30").Value = "  -6.17%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.96"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'134.09"
$ws.Range("E33").Value = "  -19.08%  "
$ws.Range("D34").Value = "'17.21"
$ws.Range("E34").Value = "  -5.33%  "
$ws.Range("D35").Value = "'0.0692"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "'2.31"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("D37").Value = "'4.32"
$ws.Range("E37").Value = "  -5.20%  "
$ws.Range("D38").Value = "'1.80"
$ws.Range("E38").Value = "  +2.09%  "
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("D40").Value = "'22.30"
$ws.Range("E40").Value = "  +23.32%  "
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").Value = "'1.926.32"
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("D44").Value = "'0.0279"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("D45").Value = "'10.13"
$ws.Range("E45").Value = "  -5.54%  "
$ws.Range("D46").Value = "'2.07"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").Value = "'2.70"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("D48").Value = "'2.87"
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").Value = "'2.554.31"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "'52.32"
$ws.Range("E50").Value = "  -2.39%  "
$ws.Range("D51").Value = "'71.98"
$ws.Range("E51").Value = "  +0.09%  "
